$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F, shifting the existing District column (old F) to G
$ws.Columns("F:F").Insert()

$ws.Range("F2").Value = "Address"
$ws.Range("G2").Value = "District"

$addresses = @{
    3 = "G H S KanagamakalapalliBagepalli"
    4 = "N G H S VenkatagirikoteChintamani"
    5 = "G H S KadalaveniGauribidanur"
    6 = "G H S MuthurSidlaghatta"
    8 = "G B H S Bagepally townBagepally"
    9 = "G H P S Poshettihalli"
    10 = "G H S YelloduGudibande"
    11 = "G H S TalagavaraCHintamani"
    12 = "G H S MelurSidlghatta"
    13 = "G H S D. PalyaGowribidanur"
    14 = "G G H S Gudibande"
    15 = "G H S KaiwaraChintamani"
    16 = "G B H S Gudibande"
    17 = "G H P S Kuduvathi"
    18 = "Gnanamandira High SchoolGowribidanur"
    20 = "Chintamani"
    21 = "G H S ChakaveluBagepalli"
    22 = "G H S Sidlaghatta"
    23 = "G H S SomenahalliGudibande"
    24 = "G G H S Chintamani"
    25 = "G H S SanthekallahalliChintamani"
    26 = "G H S MylandlahalliChintamani"
    27 = "M G H P S MurugamaleChintamani"
    28 = "G H S ChinnasandraChintamani"
    29 = "G H S BychapuraGowribidanur"
    30 = "G H S M NallagutalahalliBagepalli"
    31 = "Govt. High School K RaguttahalliChintamani"
    32 = "S S S S High School ThondebhaviGowribidanur"
    33 = "G H P S Haristala"
    34 = "Chandana High School ChamdanadoorGauribidanur"
    35 = "G H P S Kuppahalli"
    36 = "G H S HuduguruGauribidanur"
    37 = "G H S AlakapuraGauribidanur"
    38 = "S M H S VidhuraswthaGouribidanur"
    39 = "Govt. High School GhantamvaripalliBagepalli"
    40 = "G H S ThmmapalliBagepalli"
    41 = "G H S Sidlghatta"
    42 = "G H S R G Halli"
    43 = "G B H S Bagepalli"
    44 = "G M H P S HossurGowribidanur"
    45 = "G H S MuragamaleChintamani"
    46 = "G H P SChelumenahalli"
    47 = "G H S Sidlaghatta"
    48 = "G H S KundalagurkiSidlagatta"
    49 = "Rural High SchoolBurudugunteChintamani"
    50 = "G H S PalicherluShidlaghatta"
    51 = "G M H P S Peresandra"
    52 = "G H S Ullodu Gudibande"
    53 = "G H S VatadahosahallyGowribidanur"
    54 = "Govt. High School Reddy Gollavarahalli"
    55 = "G H S CheemangalaSidlaghatta taluk"
    56 = "G H S LakshmidevanakoteChintamani"
    57 = "G H P S Koothanahally"
    58 = "G H S BeechaganahalliGudibande"
    59 = "G H S MachahallyGudibande"
    60 = "G J C Gudibande"
    61 = "G H S Melya Gowribidanur"
    62 = "G H S MittemariBagepalli"
    63 = "G H S DyavappanagudiShidlghatta"
    64 = "G B H SGudibande"
    65 = "S S S S H S ThondebhaviGowribidanuru"
    66 = "G H P S SonnashettihalliChintamani"
    67 = "G H S NagergereGauribidanur"
    68 = "G H S GanjigunteShidlghatta"
    69 = "G U U H P S Taiba NagarSidlaghatta"
}

foreach ($row in $addresses.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $addresses[$row]
}
